# ST-785: Rename sheet for packaging.
# The workbook's single sheet is titled "MSL Packaging"; rename it to
# "Packaging" (workbook.xml <sheet name="..."/> attribute).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Packaging"
